$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: update Hydrogen / Iron & steel value
$ws.Range("B3").Value = 812211.0323929896

# D3: clear the Hydrogen / Non-metallic minerals value (no longer reported)
$ws.Range("D3").ClearContents()

# C4: update Methanol / Chemicals value
$ws.Range("C4").Value = 26.94949782463479

# C5: update Ammonia / Chemicals value
$ws.Range("C5").Value = 475.3446091951566

# Row 7: rename "Other" to "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 520.8128278762962

# Row 8: new "Other" row, copying row 7's formatting (border/bold/alignment)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 182.2548755457354
